# Optimising the framework with merging import statements and Code
#
# The "Passwords" label in A2 is replaced with a single blank space, and
# the sheet's active selection is moved to M12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 previously held "Passwords" - change it to a single space " ".
$ws.Range("A2").Value = " "

# Move the active cell / selection to M12 (previously E7:E8 / E8).
$ws.Range("M12").Select()
